$wb = $excel.ActiveWorkbook

# Rename Sheet2 to William
$ws = $wb.Worksheets.Item("Sheet2")
$ws.Name = "William"

# Fix the E3 formula (remove stray leading space) - value stays the same
$ws.Range("E3").Formula = "=E2+1"

# Add new "Log(+n)" labels in column F for several rows
$ws.Range("F5").Value = "Log(+3)"
$ws.Range("F6").Value = "Log(+1)"
$ws.Range("F11").Value = "Log(+1)"
$ws.Range("F16").Value = "Log(+1)"
$ws.Range("F19").Value = "Log(+1)"
$ws.Range("F22").Value = "Log(+2)"
$ws.Range("F26").Value = "Log(+2)"
$ws.Range("F32").Value = "Log(+1)"

# Update the selection shown in the sheet view
$ws.Range("F33").Select()
